$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.552.17'
Set-TextValue 'E2' '  +2.80%  '
Set-TextValue 'D3' '2.690.97'
Set-TextValue 'E3' '  +1.36%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '522.27'
Set-TextValue 'E5' '  +1.87%  '
Set-TextValue 'D6' '146.54'
Set-TextValue 'E6' '  +1.75%  '
Set-TextValue 'D7' '0.996'
Set-TextValue 'E7' '  -0.07%  '
Set-TextValue 'D8' '0.579'
Set-TextValue 'E8' '  +1.65%  '
Set-TextValue 'D9' '2.708.68'
Set-TextValue 'E9' '  +1.95%  '
Set-TextValue 'E10' '  +2.56%  '
Set-TextValue 'D11' '0.106'
Set-TextValue 'E11' '  +0.29%  '
Set-TextValue 'E12' '  +1.43%  '
Set-TextValue 'E13' '  +1.67%  '
Set-TextValue 'D14' '3.155.37'
Set-TextValue 'E14' '  +1.38%  '
Set-TextValue 'D15' '60.465.49'
Set-TextValue 'E15' '  +2.67%  '
Set-TextValue 'D16' '21.35'
Set-TextValue 'E16' '  +1.57%  '
Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.0000139'
Set-TextValue 'E17' '  +1.75%  '
Set-TextValue 'B18' 'WrappedEther'
Set-TextValue 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '2.710.04'
Set-TextValue 'E18' '  +2.05%  '
Set-TextValue 'D19' '351.56'
Set-TextValue 'E19' '  +2.43%  '
Set-TextValue 'E20' '  +0.36%  '
Set-TextValue 'D21' '10.57'
Set-TextValue 'E21' '  +2.13%  '
Set-TextValue 'E22' '  +3.99%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  -0.04%  '
Set-TextValue 'D24' '63.03'
Set-TextValue 'E24' '  +3.50%  '
Set-TextValue 'D25' '0.423'
Set-TextValue 'E25' '  +0.88%  '
Set-TextValue 'D26' '0.169'
Set-TextValue 'E26' '  +5.22%  '
Set-TextValue 'D27' '0.994'
Set-TextValue 'E27' '  -0.20%  '
Set-TextValue 'D28' '0.0₃0817'
Set-TextValue 'E28' '  +1.60%  '
Set-TextValue 'D29' '7.27'
Set-TextValue 'E29' '  +2.37%  '
Set-TextValue 'D30' '6.92'
Set-TextValue 'E30' '  +8.04%  '
Set-TextValue 'E31' '  +0.10%  '
Set-TextValue 'D32' '1.60'
Set-TextValue 'E32' '  +1.63%  '
Set-TextValue 'D33' '19.09'
Set-TextValue 'E33' '  +1.10%  '
Set-TextValue 'D34' '148.71'
Set-TextValue 'E34' '  -0.32%  '
Set-TextValue 'D35' '4.36'
Set-TextValue 'E35' '  +8.06%  '
Set-TextValue 'D36' '0.953'
Set-TextValue 'E36' '  -5.15%  '
Set-TextValue 'E37' '  +7.53%  '
Set-TextValue 'E38' '  +11.67%  '
Set-TextValue 'D39' '0.877'
Set-TextValue 'E39' '  +3.01%  '
Set-TextValue 'D40' '36.83'
Set-TextValue 'E40' '  +1.01%  '
Set-TextValue 'E41' '  +0.72%  '
Set-TextValue 'D42' '283.47'
Set-TextValue 'E42' '  +1.11%  '
Set-TextValue 'B43' 'Mantle'
Set-TextValue 'C43' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D43' '0.613'
Set-TextValue 'E43' '  -0.04%  '
Set-TextValue 'B44' 'Stellar'
Set-TextValue 'C44' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D44' '0.0991'
Set-TextValue 'E44' '  +0.85%  '
Set-TextValue 'B45' 'EnergySwap'
Set-TextValue 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '20.02'
Set-TextValue 'E45' '  +3.25%  '
Set-TextValue 'B46' 'FirstDigitalUSD'
Set-TextValue 'C46' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D46' '0.994'
Set-TextValue 'E46' '  -0.43%  '
Set-TextValue 'D47' '2.132.24'
Set-TextValue 'E47' '  +7.69%  '
Set-TextValue 'D48' '0.0542'
Set-TextValue 'E48' '  +2.02%  '
Set-TextValue 'D49' '4.87'
Set-TextValue 'E49' '  +3.96%  '
Set-TextValue 'E50' '  +3.06%  '
Set-TextValue 'D51' '10.46'
Set-TextValue 'E51' '  +1.88%  '
